$wb = $excel.ActiveWorkbook

$wsInput = $wb.Worksheets.Item("ProductLoanInput")
$wsOutput = $wb.Worksheets.Item("ProductLoanOutput")

# Update product name (B1) and short name (B2) on the input sheet
$wsInput.Range("B1").Value = "4352-MS-Simple-Group-Loan-Product-Loanproduct-1st"
$wsInput.Range("B2").Value = "435l"

# Mirror the product name change on the output sheet
$wsOutput.Range("B1").Value = "4352-MS-Simple-Group-Loan-Product-Loanproduct-1st"

# Update the selected cell on the input sheet
$wsInput.Activate()
$wsInput.Range("B3").Select()
